# DeveloperGuide: update section of UndoRedoStack to UndoRedoCareTaker
#
# 1) The auto-updating "datetimeFigureOut" date placeholders (cached display
#    text "7/20/17") get refreshed to "4/4/2018" on the slide master and on
#    every slide layout (the notes master exposes the same placeholder too,
#    but its shape collection is handled separately below).
# 2) Three now-obsolete shapes that described the old "UndoRedoStack" design
#    (a rectangle, its connector arrow, and its "1" annotation textbox) are
#    removed from the one diagram slide.

$p = $ppt.ActivePresentation

function Update-DatePlaceholder($shapesOwner, $newText) {
    $shapes = $shapesOwner.Shapes
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $shp = $shapes.Item($i)
        if ($shp.Name -like "Date Placeholder*" -and $shp.HasTextFrame) {
            if ($shp.TextFrame.HasText) {
                $shp.TextFrame.TextRange.Text = $newText
            }
        }
    }
}

# --- 1. Refresh the cached "datetimeFigureOut" field text -----------------

# Slide master
Update-DatePlaceholder $p.SlideMaster "4/4/2018"

# Every slide layout under the slide master
$layouts = $p.SlideMaster.CustomLayouts
for ($li = 1; $li -le $layouts.Count; $li++) {
    Update-DatePlaceholder $layouts.Item($li) "4/4/2018"
}

# --- 2. Remove the obsolete UndoRedoStack shapes on the diagram slide -----

function Remove-ShapeById($slide, $id) {
    for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
        $shp = $slide.Shapes.Item($i)
        if ($shp.Id -eq $id) {
            $shp.Delete()
            return
        }
    }
}

$s = $p.Slides.Item(1)

# "Rectangle 62" (UndoRedo / Stack box)
Remove-ShapeById $s 59
# "Straight Arrow Connector 57" (arrow pointing at the box)
Remove-ShapeById $s 61
# "TextBox 62" (the "1" multiplicity label)
Remove-ShapeById $s 63
